# NCATS Study multifilter testcases61to70
# Update the "startup" worksheet: replace the query text in B2 with the
# revised Cypher query (Program-level filter instead of study-subject
# filter), and remove the now-unused row 4 (the old "Beagle" breed
# filter queries), which also drops the now-unreferenced shared strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newQuery = @'
MATCH (ss:study_subject)
Match (s:study)
MATCH (ss)<-[:sample_of_study_subject]-(sp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
WITH ss, collect(DISTINCT sp.sample_id) AS samples, 
collect(DISTINCT lp.laboratory_procedure_id) AS lab_procedures, 
collect(DISTINCT f) AS files
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
OPTIONAL MATCH (ss)<-[:program_of_institution]-(p)
OPTIONAL MATCH (p)<-[:of_arm]-(a)
MATCH (ss)<-[:sf_of_study_subject]-(sf)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
MATCH (ss)<-[:demographic_of_study_subject]-(demo)
RETURN DISTINCT
       coalesce (p.program_acronym, '')as `Program Code`,
       coalesce( p.program_id , '')as `Program ID`,
       coalesce (p.program_name , '' )as `Program Name`,
       coalesce(p.start_date, '') as `Start Date`,
       coalesce (p.end_date, '') as `End Date`,
       coalesce(p.pubmed_id, '') as `PubMed ID`,
       count(distinct s) As `Number of Arms`,
       count(distinct ss) as `Associated Cases`
       order By `Program Code`
'@

# Update the query stored in B2 (trim the trailing newline the here-string adds)
$ws.Range("B2").Value = $newQuery.TrimEnd("`r", "`n")

# Row 4 (the old breed/"Beagle" filter queries) is no longer used - delete it
# entirely so the remaining rows collapse and the dimension shrinks to A1:E2.
$ws.Rows(4).Delete()

# Row 2 grew a little shorter once the old query text left the sheet -
# re-apply the row height used for the refreshed query text.
$ws.Rows(2).RowHeight = 259.5

# Minor column-width touch ups that came along with the content edit.
$ws.Columns(1).ColumnWidth = 12.9555
$ws.Columns(2).ColumnWidth = 74.7609
$ws.Columns(4).ColumnWidth = 40.4352
$ws.Columns(5).ColumnWidth = 42.7609
